$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9179683923721313
$ws.Range("B1").Value = 0.778282105922699
$ws.Range("C1").Value = 3.667492866516113
$ws.Range("D1").Value = 2.85458779335022
$ws.Range("E1").Value = 1.303042769432068
